# Update scripts with new TPM values (Fgf8-Fgfr1, LR-pairs_lrc2p, YoungD4)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Sending cluster" column relabelled as part of the regenerated export
$ws.Range("A2").Value = "Resolving-Mac"
$ws.Range("A3").Value = "Resolving-Mac"
$ws.Range("A4").Value = "Resolving-Mac"
$ws.Range("A5").Value = "Resolving-Mac"
$ws.Range("A6").Value = "Resolving-Mac"

# Row 2 (target: ECs) - updated TPM-derived values
$ws.Range("G2").Value = 0.02530666666666667
$ws.Range("H2").Value = 0.07592
$ws.Range("M2").Value = 13.67700833333333
$ws.Range("N2").Value = 41.031025
$ws.Range("O2").Value = 0.124413831206147
$ws.Range("P2").Value = 0.124413831206147
$ws.Range("Q2").Value = 0.3461194908888889
$ws.Range("R2").Value = 3.115075418
$ws.Range("S2").Value = 0.124413831206147
$ws.Range("T2").Value = 0.124413831206147

# Row 3 (target: FAPs) - updated TPM-derived values
$ws.Range("G3").Value = 0.02530666666666667
$ws.Range("H3").Value = 0.07592
$ws.Range("M3").Value = 74.64939600000001
$ws.Range("N3").Value = 223.948188
$ws.Range("O3").Value = 0.679053278848249
$ws.Range("P3").Value = 0.6790532788482488
$ws.Range("Q3").Value = 1.88912738144
$ws.Range("R3").Value = 17.00214643296
$ws.Range("S3").Value = 0.679053278848249
$ws.Range("T3").Value = 0.6790532788482488

# Row 4 (target: Inflammatory-Mac) - updated TPM-derived values
$ws.Range("G4").Value = 0.02530666666666667
$ws.Range("H4").Value = 0.07592
$ws.Range("M4").Value = 1.629335666666667
$ws.Range("N4").Value = 4.888007
$ws.Range("O4").Value = 0.01482136207497777
$ws.Range("P4").Value = 0.01482136207497777
$ws.Range("Q4").Value = 0.04123305460444445
$ws.Range("R4").Value = 0.37109749144
$ws.Range("S4").Value = 0.01482136207497777
$ws.Range("T4").Value = 0.01482136207497777

# Row 5 (target: MuSCs) - updated TPM-derived values
$ws.Range("G5").Value = 0.02530666666666667
$ws.Range("H5").Value = 0.07592
$ws.Range("M5").Value = 19.17462033333333
$ws.Range("N5").Value = 57.523861
$ws.Range("O5").Value = 0.174423230537864
$ws.Range("P5").Value = 0.174423230537864
$ws.Range("Q5").Value = 0.4852457252355556
$ws.Range("R5").Value = 4.367211527119999
$ws.Range("S5").Value = 0.174423230537864
$ws.Range("T5").Value = 0.174423230537864

# Row 6 (target: Resolving-Mac) - updated TPM-derived values
$ws.Range("G6").Value = 0.02530666666666667
$ws.Range("H6").Value = 0.07592
$ws.Range("M6").Value = 0.801214
$ws.Range("N6").Value = 2.403642
$ws.Range("O6").Value = 0.007288297332762355
$ws.Range("P6").Value = 0.007288297332762355
$ws.Range("Q6").Value = 0.02027605562666667
$ws.Range("R6").Value = 0.18248450064
$ws.Range("S6").Value = 0.007288297332762355
$ws.Range("T6").Value = 0.007288297332762355
